$wb = $excel.ActiveWorkbook

# Overview sheet: Correspond Handback DateTime (G4) for cc5f8fb9 row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-17 10:42:35"

# zh-cn sheet: Correspond Handoff Datetime (H4) and Correspond Handback DateTime (K4) for cc5f8fb9 row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-17 10:42:30"
$wsZhCn.Range("K4").Value = "2016-08-17 10:42:48"

# de-de sheet: Correspond Handoff Datetime (H4) and Correspond Handback DateTime (K4) for cc5f8fb9 row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-17 10:42:35"
$wsDeDe.Range("K4").Value = "2016-08-17 10:42:55"
